# Insert two new weekly records (Asterix "1a (guarda)" and a new Rosara
# "1a (cosecha)" entry) right after the existing row for Ñuble / Papa,
# pushing the historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 679-680; everything that used to start at row 679
# shifts down to row 681 onward.
$ws.Range("A679:A680").EntireRow.Insert()

# New row 679: Asterix, "1a (guarda)"
$ws.Range("A679").Value = 7
$ws.Range("B679").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C679").Value = "Ñuble"
$ws.Range("D679").Value = 45147
$ws.Range("E679").Value = 16
$ws.Range("F679").Value = 100114001
$ws.Range("G679").Value = "Papa"
$ws.Range("H679").Value = "Asterix"
$ws.Range("I679").Value = "1a (guarda)"
$ws.Range("J679").Value = 100
$ws.Range("K679").Value = 19000
$ws.Range("L679").Value = 19000
$ws.Range("M679").Value = 19000
$ws.Range("N679").Value = '$/saco 25 kilos'
$ws.Range("O679").Value = "Región de Los Lagos"
$ws.Range("P679").Value = 760
$ws.Range("Q679").Value = 25
$ws.Range("R679").Value = "Hortaliza"

# New row 680: Rosara, "1a (cosecha)"
$ws.Range("A680").Value = 7
$ws.Range("B680").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C680").Value = "Ñuble"
$ws.Range("D680").Value = 45147
$ws.Range("E680").Value = 16
$ws.Range("F680").Value = 100114001
$ws.Range("G680").Value = "Papa"
$ws.Range("H680").Value = "Rosara"
$ws.Range("I680").Value = "1a (cosecha)"
$ws.Range("J680").Value = 50
$ws.Range("K680").Value = 17000
$ws.Range("L680").Value = 17000
$ws.Range("M680").Value = 17000
$ws.Range("N680").Value = '$/saco 25 kilos'
$ws.Range("O680").Value = "Región del Maule"
$ws.Range("P680").Value = 680
$ws.Range("Q680").Value = 25
$ws.Range("R680").Value = "Hortaliza"
